# Generate Report for Handback
# Updates timestamps / status in the handback status report workbook.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 06:16:57"
$wsOverview.Range("G5").Value = "2016-08-19 06:16:57"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-19 06:16:51"
$wsZhCn.Range("H5").Value = "2016-08-19 06:16:51"
$wsZhCn.Range("K2").Value = "2016-08-19 06:17:13"
$wsZhCn.Range("K5").Value = "2016-08-19 06:17:13"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-19 06:16:57"
$wsDeDe.Range("H5").Value = "2016-08-19 06:16:57"
$wsDeDe.Range("K2").Value = "2016-08-19 06:17:20"
$wsDeDe.Range("K5").Value = "2016-08-19 06:17:20"
